$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a handful of Supplier Stock 1 (column J) quantities that changed
# independently of the row removal below.
$ws.Cells.Item(2, 10).Value = 1903578
$ws.Cells.Item(7, 10).Value = 43490
$ws.Cells.Item(17, 10).Value = 4399293

# Revert to the old/smaller basket size: remove the extra "B1" / "N-5"
# line item (row 23) that had been added, shifting everything below it
# back up by one row.
$ws.Rows(23).Delete()

# The "report created" time label (text, next to the printed: NOW() cell)
# needs to be refreshed to match the new basket number / recalculation.
$ws.Cells.Item(24, 6).Value = "17:43"
